$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the clinical-variables list in column D (page4): insert the new
# "Pembro_status" variable, re-sequence, and append the several new
# clinical variables that came with it (Sample_type, Body_part,
# Body_part_abbrev, Metastasis_brain_status, predictive_biomarker,
# Best_overall_response_group/_detailed, PFS/OS month groupings,
# Impact_TMB_score_group10, STK11, KEAP1). Metastasis_multiple and
# Best_overall_response_abbrev are dropped from the list.
$ws.Range("D2").Value  = "Pembro_status"
$ws.Range("D3").Value  = "Sex"
$ws.Range("D4").Value  = "Age_at_diagnosis_abbrev"
$ws.Range("D5").Value  = "Type"
$ws.Range("D6").Value  = "Sample_type"
$ws.Range("D7").Value  = "Body_part"
$ws.Range("D8").Value  = "Body_part_abbrev"
$ws.Range("D9").Value  = "Metastasis_brain"
$ws.Range("D10").Value = "Metastasis_brain_status"
$ws.Range("D11").Value = "Smoking_status_abbrev"
$ws.Range("D12").Value = "Patient_history_of_cancer_abbrev"
$ws.Range("D13").Value = "Family_history_of_cancer_abbrev"
$ws.Range("D14").Value = "predictive_biomarker"
$ws.Range("D15").Value = "Best_overall_response_group"
$ws.Range("D16").Value = "Best_overall_response_detailed"
$ws.Range("D17").Value = "PFS_months_group_median"
$ws.Range("D18").Value = "PFS_months_group_quartile"
$ws.Range("D19").Value = "OS_months_group_median"
$ws.Range("D20").Value = "OS_months_group_quartile"
$ws.Range("D21").Value = "PDL1_percent_score_group"
$ws.Range("D22").Value = "Impact_TMB_score_group10"
$ws.Range("D23").Value = "STK11"
$ws.Range("D24").Value = "KEAP1"

# The whole of column B previously carried a stray "applyFill" cell style
# (a no-op fill of "none", left over in the workbook) - clear it so the
# column goes back to the default/unstyled look, matching the cleaned-up
# workbook.
$ws.Columns("B").ClearFormats()
